$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aVals = @(45918,45918.01041666666,45918.02083333334,45918.03125,45918.04166666666,45918.05208333334,45918.0625,45918.07291666666,45918.08333333334,45918.09375,45918.10416666666,45918.11458333334,45918.125,45918.13541666666,45918.14583333334,45918.15625,45918.16666666666,45918.17708333334,45918.1875,45918.19791666666,45918.20833333334,45918.21875,45918.22916666666,45918.23958333334,45918.25,45918.26041666666,45918.27083333334,45918.28125,45918.29166666666,45918.30208333334,45918.3125,45918.32291666666,45918.33333333334,45918.34375,45918.35416666666,45918.36458333334,45918.375,45918.38541666666,45918.39583333334,45918.40625,45918.41666666666,45918.42708333334,45918.4375,45918.44791666666,45918.45833333334,45918.46875,45918.47916666666,45918.48958333334,45918.5,45918.51041666666,45918.52083333334,45918.53125,45918.54166666666,45918.55208333334,45918.5625,45918.57291666666,45918.58333333334,45918.59375,45918.60416666666,45918.61458333334,45918.625,45918.63541666666,45918.64583333334,45918.65625,45918.66666666666,45918.67708333334,45918.6875,45918.69791666666,45918.70833333334,45918.71875,45918.72916666666,45918.73958333334,45918.75,45918.76041666666,45918.77083333334,45918.78125,45918.79166666666,45918.80208333334,45918.8125,45918.82291666666,45918.83333333334,45918.84375,45918.85416666666,45918.86458333334,45918.875,45918.88541666666,45918.89583333334,45918.90625,45918.91666666666,45918.92708333334,45918.9375,45918.94791666666,45918.95833333334,45918.96875,45918.97916666666,45918.98958333334,45919,45919.01041666666,45919.02083333334,45919.03125,45919.04166666666,45919.05208333334,45919.0625,45919.07291666666,45919.08333333334,45919.09375,45919.10416666666,45919.11458333334,45919.125,45919.13541666666,45919.14583333334,45919.15625,45919.16666666666,45919.17708333334,45919.1875,45919.19791666666,45919.20833333334,45919.21875,45919.22916666666,45919.23958333334,45919.25,45919.26041666666,45919.27083333334,45919.28125,45919.29166666666,45919.30208333334,45919.3125,45919.32291666666,45919.33333333334,45919.34375,45919.35416666666,45919.36458333334,45919.375,45919.38541666666,45919.39583333334,45919.40625,45919.41666666666,45919.42708333334,45919.4375,45919.44791666666,45919.45833333334,45919.46875,45919.47916666666,45919.48958333334,45919.5,45919.51041666666,45919.52083333334,45919.53125,45919.54166666666,45919.55208333334,45919.5625,45919.57291666666,45919.58333333334,45919.59375,45919.60416666666,45919.61458333334,45919.625,45919.63541666666,45919.64583333334,45919.65625,45919.66666666666,45919.67708333334,45919.6875,45919.69791666666,45919.70833333334,45919.71875,45919.72916666666,45919.73958333334,45919.75,45919.76041666666,45919.77083333334,45919.78125,45919.79166666666,45919.80208333334,45919.8125,45919.82291666666,45919.83333333334,45919.84375,45919.85416666666,45919.86458333334,45919.875,45919.88541666666,45919.89583333334,45919.90625,45919.91666666666,45919.92708333334,45919.9375,45919.94791666666,45919.95833333334,45919.96875,45919.97916666666,45919.98958333334)
$bVals = @(5184,5139,5137,5095,4961,5014,4980,4962,4860,4926,4970,4926,4970,4936,5016,4971,5096,5138,5151,5247,5446,5551,5678,5764,5974,6151,6227,6263,6255,6172,6178,6171,5883,5818,5783,5616,5324,5195,5135,5055,4914,4844,4810,4752,4690,4688,4727,4691,4725,4803,4827,4828,4886,4837,4901,4910,5027,4980,5031,5112,5259,5284,5296,5374,5464,5561,5753,5853,6090,6220,6326,6500,6550,6674,6838,6938,7117,7117,7113,7033,6950,6868,6720,6604,6332,6214,6033,5935,5699,5548,5479,5402,5291,5246,5187,5150,5053,5063,5051,5031,4981,5038,4972,4985,4933,4949,4906,4922,4942,5011,4974,4991,5117,5198,5168,5262,5484,5587,5640,5709,5941,6121,6129,6160,6230,6284,6235,6082,5983,5905,5815,5680,5577,5473,5343,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
$dVals = @("18.09.20251","18.09.20252","18.09.20253","18.09.20254","18.09.20255","18.09.20256","18.09.20257","18.09.20258","18.09.20259","18.09.202510","18.09.202511","18.09.202512","18.09.202513","18.09.202514","18.09.202515","18.09.202516","18.09.202517","18.09.202518","18.09.202519","18.09.202520","18.09.202521","18.09.202522","18.09.202523","18.09.202524","18.09.202525","18.09.202526","18.09.202527","18.09.202528","18.09.202529","18.09.202530","18.09.202531","18.09.202532","18.09.202533","18.09.202534","18.09.202535","18.09.202536","18.09.202537","18.09.202538","18.09.202539","18.09.202540","18.09.202541","18.09.202542","18.09.202543","18.09.202544","18.09.202545","18.09.202546","18.09.202547","18.09.202548","18.09.202549","18.09.202550","18.09.202551","18.09.202552","18.09.202553","18.09.202554","18.09.202555","18.09.202556","18.09.202557","18.09.202558","18.09.202559","18.09.202560","18.09.202561","18.09.202562","18.09.202563","18.09.202564","18.09.202565","18.09.202566","18.09.202567","18.09.202568","18.09.202569","18.09.202570","18.09.202571","18.09.202572","18.09.202573","18.09.202574","18.09.202575","18.09.202576","18.09.202577","18.09.202578","18.09.202579","18.09.202580","18.09.202581","18.09.202582","18.09.202583","18.09.202584","18.09.202585","18.09.202586","18.09.202587","18.09.202588","18.09.202589","18.09.202590","18.09.202591","18.09.202592","18.09.202593","18.09.202594","18.09.202595","18.09.202596","19.09.20251","19.09.20252","19.09.20253","19.09.20254","19.09.20255","19.09.20256","19.09.20257","19.09.20258","19.09.20259","19.09.202510","19.09.202511","19.09.202512","19.09.202513","19.09.202514","19.09.202515","19.09.202516","19.09.202517","19.09.202518","19.09.202519","19.09.202520","19.09.202521","19.09.202522","19.09.202523","19.09.202524","19.09.202525","19.09.202526","19.09.202527","19.09.202528","19.09.202529","19.09.202530","19.09.202531","19.09.202532","19.09.202533","19.09.202534","19.09.202535","19.09.202536","19.09.202537","19.09.202538","19.09.202539","19.09.202540","19.09.202541","19.09.202542","19.09.202543","19.09.202544","19.09.202545","19.09.202546","19.09.202547","19.09.202548","19.09.202549","19.09.202550","19.09.202551","19.09.202552","19.09.202553","19.09.202554","19.09.202555","19.09.202556","19.09.202557","19.09.202558","19.09.202559","19.09.202560","19.09.202561","19.09.202562","19.09.202563","19.09.202564","19.09.202565","19.09.202566","19.09.202567","19.09.202568","19.09.202569","19.09.202570","19.09.202571","19.09.202572","19.09.202573","19.09.202574","19.09.202575","19.09.202576","19.09.202577","19.09.202578","19.09.202579","19.09.202580","19.09.202581","19.09.202582","19.09.202583","19.09.202584","19.09.202585","19.09.202586","19.09.202587","19.09.202588","19.09.202589","19.09.202590","19.09.202591","19.09.202592","19.09.202593","19.09.202594","19.09.202595","19.09.202596")

for ($i = 0; $i -lt 192; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aVals[$i]
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
}
